$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pedir Wurth")

# Copy plain border style (C3 = style index "1") onto the new cells first, so that
# typing into them afterwards reuses the existing style entries instead of
# allocating brand-new ones.
$ws.Range("C3").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C12:C13").PasteSpecial(-4122)
$ws.Range("D12:E13").PasteSpecial(-4122)

$ws.Cells.Item(12,1).Value = "'U9"
$ws.Cells.Item(12,2).Value = "'Opto Triac 400VDRM 10mA"
$ws.Cells.Item(12,3).Value = 10
$ws.Cells.Item(12,4).Value = "'Wurth Electronics"
$ws.Cells.Item(12,5).Value = "'14230226011"

$ws.Cells.Item(13,1).Value = "'U6, U7, U8"
$ws.Cells.Item(13,2).Value = "'Opto Coupler"
$ws.Cells.Item(13,3).Value = 30
$ws.Cells.Item(13,4).Value = "'Wurth Electronics"
$ws.Cells.Item(13,5).Value = "'140816140410"

Write-Host "done"
